# Auto-generated edit script applying the cryptos.xlsx price/volume update
# (commit: "Updated cryptos list on Thu Jan 18 18:43:35 UTC 2024 with GitHub Actions")
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D2").Value = "41.538.11"
$ws.Range("E2").Value = "  -2.13%  "
$ws.Range("D3").Value = "2.465.12"
$ws.Range("E3").Value = "  -2.40%  "
$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").Value = "'310.81"
$ws.Range("E5").Value = "  +0.70%  "
$ws.Range("D6").Value = "'95.02"
$ws.Range("E6").Value = "  -4.29%  "
$ws.Range("D7").Value = "'0.552"
$ws.Range("E7").Value = "  -2.92%  "
$ws.Range("E8").Value = "  +0.02%  "
$ws.Range("D9").Value = "'0.509"
$ws.Range("E9").Value = "  -3.49%  "
$ws.Range("D10").Value = "'33.87"
$ws.Range("E10").Value = "  -5.33%  "
$ws.Range("D11").Value = "'0.0784"
$ws.Range("E11").Value = "  -2.41%  "
$ws.Range("E12").Value = "  -0.08%  "
$ws.Range("D13").Value = "'6.99"
$ws.Range("E13").Value = "  -4.42%  "
$ws.Range("D14").Value = "2.841.85"
$ws.Range("E14").Value = "  -2.48%  "
$ws.Range("D15").Value = "2.455.32"
$ws.Range("E15").Value = "  -2.46%  "
$ws.Range("D16").Value = "'14.58"
$ws.Range("E16").Value = "  -7.48%  "
$ws.Range("D17").Value = "'0.790"
$ws.Range("E17").Value = "  -3.54%  "
$ws.Range("D18").Value = "41.524.64"
$ws.Range("E18").Value = "  -2.12%  "
$ws.Range("D19").Value = "'6.39"
$ws.Range("E19").Value = "  -5.95%  "
$ws.Range("D20").Value = "0.0₃0919"
$ws.Range("E20").Value = "  -3.56%  "
$ws.Range("D21").Value = "'11.56"
$ws.Range("E21").Value = "  -4.99%  "
$ws.Range("D22").Value = "'69.69"
$ws.Range("E22").Value = "  +0.98%  "
$ws.Range("D23").Value = "'237.19"
$ws.Range("E23").Value = "  -2.06%  "
$ws.Range("D24").Value = "'2.78"
$ws.Range("E24").Value = "  -3.44%  "
$ws.Range("D25").Value = "'1.94"
$ws.Range("E25").Value = "  -4.50%  "
$ws.Range("E26").Value = "  +0.07%  "
$ws.Range("D27").Value = "'24.76"
$ws.Range("E27").Value = "  -4.36%  "
$ws.Range("E28").Value = "  -5.12%  "
$ws.Range("D29").Value = "'9.75"
$ws.Range("E29").Value = "  -3.45%  "
$ws.Range("D30").Value = "'36.35"
$ws.Range("E30").Value = "  -7.39%  "
$ws.Range("D31").Value = "'153.08"
$ws.Range("E31").Value = "  -1.73%  "
$ws.Range("D32").Value = "'5.64"
$ws.Range("E32").Value = "  -1.55%  "
$ws.Range("E33").Value = "  -0.24%  "
$ws.Range("B34").Value = "Hedera"
$ws.Range("C34").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D34").Value = "'0.0756"
$ws.Range("E34").Value = "  -4.58%  "
$ws.Range("B35").Value = "ApeXProtocol"
$ws.Range("C35").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D35").Value = "'2.54"
$ws.Range("E35").Value = "  -8.67%  "
$ws.Range("D36").Value = "'3.02"
$ws.Range("E36").Value = "  -4.64%  "
$ws.Range("D37").Value = "'17.23"
$ws.Range("E37").Value = "  -5.19%  "
$ws.Range("D38").Value = "'1.88"
$ws.Range("E38").Value = "  -6.76%  "
$ws.Range("E39").Value = "  -5.15%  "
$ws.Range("E40").Value = "  -3.19%  "
$ws.Range("D41").Value = "'4.03"
$ws.Range("E41").Value = "  -6.43%  "
$ws.Range("D42").Value = "'21.47"
$ws.Range("E42").Value = "  -2.59%  "
$ws.Range("D44").Value = "1.983.51"
$ws.Range("E44").Value = "  +0.83%  "
$ws.Range("E45").Value = "  -4.18%  "
$ws.Range("D46").Value = "'3.07"
$ws.Range("E46").Value = "  -7.12%  "
$ws.Range("D47").Value = "'8.71"
$ws.Range("E47").Value = "  -1.99%  "
$ws.Range("D48").Value = "2.697.39"
$ws.Range("E48").Value = "  -2.61%  "
$ws.Range("D49").Value = "'76.82"
$ws.Range("E49").Value = "  -5.18%  "
$ws.Range("B50").Value = "Aave"
$ws.Range("C50").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D50").Value = "'97.89"
$ws.Range("E50").Value = "  -3.20%  "
$ws.Range("B51").Value = "ordi"
$ws.Range("C51").Value = "https://coinranking.com/coin/j7-7vPrOi+ordi-ordi"
$ws.Range("D51").Value = "'69.64"
$ws.Range("E51").Value = "  -3.47%  "
